$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -8
$ws.Range("F3").Value = -7
$ws.Range("F7").Value = -5
$ws.Range("F8").Value = 7
$ws.Range("F10").Value = 5
$ws.Range("F13").Value = 5
$ws.Range("F14").Value = -9
$ws.Range("F18").Value = -4
$ws.Range("F19").Value = -7
$ws.Range("F20").Value = -2
$ws.Range("F21").Value = -4
$ws.Range("F22").Value = -3
$ws.Range("F23").Value = -2
$ws.Range("F24").Value = -1
$ws.Range("F25").Value = 3
$ws.Range("F26").Value = 4
$ws.Range("F27").Value = 1
$ws.Range("F28").Value = 4
$ws.Range("F29").Value = -5
$ws.Range("F30").Value = -6
$ws.Range("F31").Value = 5
$ws.Range("F33").Value = 3
$ws.Range("F34").Value = 7
$ws.Range("F38").Value = 2
$ws.Range("F39").Value = -2
$ws.Range("F41").Value = -1
$ws.Range("F42").Value = -4
$ws.Range("F43").Value = -5
$ws.Range("F44").Value = -2
$ws.Range("F45").Value = 3
